$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric are stored as text,
# matching the original inline-string formatting (e.g. "237.20", "1.000").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.773.31"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.753.13"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.20"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5056"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.61"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2622"
$ws.Range("E9").Value = "  +8.41%  "
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.748.67"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06947"
$ws.Range("E12").Value = "  +3.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.45"
$ws.Range("E13").Value = "  +3.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6029"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.53"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.451"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.811.89"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.65"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006765"
$ws.Range("E21").Value = "  +6.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.975.21"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.051"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.197"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.170"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.90"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.464"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.08"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.796"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.28"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08273"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.697"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.385"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04362"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6003"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.698"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("E40").Value = "  -5.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01546"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.34"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7471"
$ws.Range("E44").Value = "  -5.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3801"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.877"
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05480"
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1076"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.941"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.17"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  -0.09%  "
